# Adds 5 race tracks to the game data:
#  - Renames Blad1/Blad2/Blad3 to Drivers / Car Upgrades / Race Time Calculation
#  - Populates the "Race Time Calculation" sheet (was Blad3, empty) with a
#    track-comparison table plus explanatory notes
#  - Makes "Race Time Calculation" the active/selected sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename + reorder worksheets (Chart1 chart-sheet keeps trailing position)
# ---------------------------------------------------------------------------
$wsDrivers = $wb.Sheets.Item("Blad1")
$wsDrivers.Name = "Drivers"

$wsUpgrades = $wb.Sheets.Item("Blad2")
$wsUpgrades.Name = "Car Upgrades"

$wsRace = $wb.Sheets.Item("Blad3")
$wsRace.Name = "Race Time Calculation"

# ---------------------------------------------------------------------------
# 2. Cell values for the new table (row 2 header / row 3 sub-header / rows
#    4-8 per-track data / row 10 explanatory notes)
# ---------------------------------------------------------------------------
$values = @{
  "B2" = "Track Name"
  "C2" = "Amount of turns"
  "D2" = "Track total length"
  "E2" = "Track max height difference"
  "F2" = "Longest straight track"

  "C3" = "Independent of turn size"
  "D3" = "In Kilometers"
  "E3" = "In Meters"
  "F3" = "In Kilometers"

  "B4" = "Baku City Circuit, Azerbaijan"
  "C4" = "20 Corners"
  "D4" = "6,003 km"
  "E4" = "26,8 m"
  "F4" = "2,2 km"

  "B5" = "Spa-Francorchamps, Belgium"
  "C5" = "20 Corners"
  "D5" = "7,004 km"
  "E5" = "102,2 m"
  "F5" = "2,4 km"

  "B6" = "Silverstone, United Kingdom"
  "C6" = "18 Corners"
  "D6" = "5,891 km"
  "E6" = "11,3 m"
  "F6" = "1,0 km"

  "B7" = "Monza, Italy"
  "C7" = "11 Corners"
  "D7" = "5,793 km"
  "E7" = "12,8 m"
  "F7" = "1,1 km"

  "B8" = "Monte Carlo, Monaco"
  "C8" = "19 Corners"
  "D8" = "3,337 km"
  "E8" = "42,0 m"
  "F8" = "0,7 km"

  "C10" = "Car's with better handling and acceleration have better performance on track's with more turns"
  "D10" = "Longer track lengths will cause longer racing times"
  "E10" = "Track's with better down force and more weight will have less problems with height difference"
  "F10" = "Track's with more speed will be faster on track's with longer straight tracks"
}

foreach ($addr in $values.Keys) {
  $wsRace.Range($addr).Value = $values[$addr]
}

# ---------------------------------------------------------------------------
# 3. Column widths / row heights
# ---------------------------------------------------------------------------
$wsRace.Range("B1:F10").ColumnWidth = 26.7
$wsRace.Range("B1:F1").RowHeight = 15.75
$wsRace.Range("B8:F8").RowHeight = 15.75
$wsRace.Range("B10:F10").RowHeight = 52.5

# ---------------------------------------------------------------------------
# 4. Fonts: bold header row, italic sub-header row, size-10 notes row
# ---------------------------------------------------------------------------
$headerRange = $wsRace.Range("B2:F2")
$headerRange.Font.Bold = $true

$subHeaderRange = $wsRace.Range("B3:F3")
$subHeaderRange.Font.Italic = $true

$notesRange = $wsRace.Range("C10:F10")
$notesRange.Font.Size = 10

# ---------------------------------------------------------------------------
# 5. Horizontal alignment (everything centered) + notes wrap/vertical-center
# ---------------------------------------------------------------------------
$wsRace.Range("B2:F8").HorizontalAlignment = -4108
$notesRange.HorizontalAlignment = -4108
$notesRange.VerticalAlignment = -4108
$notesRange.WrapText = $true

# ---------------------------------------------------------------------------
# 6. Borders - outer box (medium) + inner column dividers / header divider
#    (thin) around B2:F8, matching the original author's manual formatting
# ---------------------------------------------------------------------------
$xlEdgeLeft = 7; $xlEdgeTop = 8; $xlEdgeBottom = 9; $xlEdgeRight = 10
$xlContinuous = 1
$xlThin = 2
$xlMedium = -4138

# left edge of the table
$colB = $wsRace.Range("B2:B8")
$colB.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$colB.Borders.Item($xlEdgeLeft).Weight = $xlMedium

# thin divider between Track Name / Amount of turns columns
$colC = $wsRace.Range("C2:C8")
$colC.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$colC.Borders.Item($xlEdgeLeft).Weight = $xlThin
$colC.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$colC.Borders.Item($xlEdgeRight).Weight = $xlThin

# thin divider between Track max height difference / In Meters columns
$colE = $wsRace.Range("E2:E8")
$colE.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$colE.Borders.Item($xlEdgeLeft).Weight = $xlThin
$colE.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$colE.Borders.Item($xlEdgeRight).Weight = $xlThin

# right edge of the table
$colF = $wsRace.Range("F2:F8")
$colF.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$colF.Borders.Item($xlEdgeRight).Weight = $xlMedium

# top edge (header row)
$headerRange.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$headerRange.Borders.Item($xlEdgeTop).Weight = $xlMedium

# thin rule under the sub-header row
$subHeaderRange.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$subHeaderRange.Borders.Item($xlEdgeBottom).Weight = $xlThin

# bottom edge (last data row)
$lastRowRange = $wsRace.Range("B8:F8")
$lastRowRange.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$lastRowRange.Borders.Item($xlEdgeBottom).Weight = $xlMedium

# ---------------------------------------------------------------------------
# 7. Activate the new sheet as the selected tab, matching the saved view
# ---------------------------------------------------------------------------
$wsRace.Activate() | Out-Null
$wsRace.Range("F21").Select() | Out-Null
